$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.001.97"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.561.43"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.27%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "207.93"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.18%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "22.11"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.781.36"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "1.557.70"
$ws.Range("E13").Value = "  +0.07%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.979.38"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "61.92"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +1.94%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "216.29"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.70%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +0.14%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.14"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.83%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.25"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -0.37%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.79"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.62"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "15.10"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("E32").Value = "  +0.00%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.16"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "1.430.36"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("E36").Value = "  +8.21%  "
$ws.Range("E37").Value = "  +2.03%  "
$ws.Range("E38").Value = "  +0.43%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.535"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.41%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.85"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.63%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.809"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("E44").Value = "  +1.46%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "64.89"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "1.696.84"
$ws.Range("E47").Value = "  +0.28%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "87.40"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  +4.85%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0518"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0962"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
